# Corrections in the guidance: append a small worked example (rows 19-22)
# below the existing calculations on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B19").Value = 3
$ws.Range("C19").Value = 20

$ws.Range("B20").Value = 106
$ws.Range("C20").Value = 26

$ws.Range("B21").Value = 805
$ws.Range("C21").Value = 77

$ws.Range("B22").Formula = "=SUM(B19:B21)"
$ws.Range("C22").Formula = "=SUM(C19:C21)"

$ws.Range("D22").Select()
